$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览 (Exhibition) -- simple numeric "want to go" (F column) bumps
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 28
$ws1.Range("F3").Value = 799
$ws1.Range("F4").Value = 38
$ws1.Range("F7").Value = 1151
$ws1.Range("F8").Value = 908
$ws1.Range("F10").Value = 718
$ws1.Range("F12").Value = 1437
$ws1.Range("F13").Value = 56
$ws1.Range("F15").Value = 1603
$ws1.Range("F17").Value = 608
$ws1.Range("F21").Value = 1079
$ws1.Range("F22").Value = 1511
$ws1.Range("F23").Value = 752
$ws1.Range("F24").Value = 616
$ws1.Range("F25").Value = 489
$ws1.Range("F30").Value = 299
$ws1.Range("F31").Value = 2422
$ws1.Range("F33").Value = 1354
$ws1.Range("F34").Value = 460
$ws1.Range("F35").Value = 66
$ws1.Range("F36").Value = 3949

# ---------------------------------------------------------------------------
# Sheet 2: 演出 (Performance)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Numeric bumps on rows that are not affected by the later row insertion.
$ws2.Range("F4").Value = 1035
$ws2.Range("F6").Value = 174
$ws2.Range("F20").Value = 256
$ws2.Range("F21").Value = 257
$ws2.Range("F23").Value = 6

# Insert two brand-new rows above the old row 29 (Rie fu concert), shifting
# everything from row 29 down onward by two rows.
$ws2.Rows.Item(29).Resize(2).Insert()

# New row 29: KAREN SOUZA concert
$ws2.Cells.Item(29, 1).Font.Bold = $true
$ws2.Cells.Item(29, 1).HorizontalAlignment = -4108
$ws2.Cells.Item(29, 1).VerticalAlignment = -4160
$ws2.Cells.Item(29, 1).Borders.LineStyle = 1
$ws2.Cells.Item(29, 1).Value = 28
$ws2.Cells.Item(29, 2).NumberFormat = "@"
$ws2.Cells.Item(29, 2).Value = "2024-05-03"
$ws2.Cells.Item(29, 2).Style = "Normal"
$ws2.Cells.Item(29, 3).Value = "上海·爵士情调女王KAREN SOUZA凯伦索萨2024演唱会"
$ws2.Cells.Item(29, 4).Value = "南京西路1376号上海商城4层 商城剧院"
$ws2.Cells.Item(29, 5).Value = "2024.05.03 19:30-05.03 21:00"
$ws2.Cells.Item(29, 6).Value = 0
$ws2.Cells.Item(29, 7).Value = 126
$ws2.Cells.Item(29, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82653"
$ws2.Cells.Item(29, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/IkBVehui1710141982443.jpeg"

# New row 30: Pianoboy concert
$ws2.Cells.Item(30, 1).Font.Bold = $true
$ws2.Cells.Item(30, 1).HorizontalAlignment = -4108
$ws2.Cells.Item(30, 1).VerticalAlignment = -4160
$ws2.Cells.Item(30, 1).Borders.LineStyle = 1
$ws2.Cells.Item(30, 1).Value = 29
$ws2.Cells.Item(30, 2).NumberFormat = "@"
$ws2.Cells.Item(30, 2).Value = "2024-05-04"
$ws2.Cells.Item(30, 2).Style = "Normal"
$ws2.Cells.Item(30, 3).Value = "上海·钢琴诗人Pianoboy高至豪流行钢琴音乐会"
$ws2.Cells.Item(30, 4).Value = "南京西路1376号上海商城4层 商城剧院"
$ws2.Cells.Item(30, 5).Value = "2024.05.04 19:30-05.04 21:00"
$ws2.Cells.Item(30, 6).Value = 0
$ws2.Cells.Item(30, 7).Value = 126
$ws2.Cells.Item(30, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82673"
$ws2.Cells.Item(30, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/MooHY44M1710149484564.jpeg"

# The three rows that used to be 29/30/31 are now 31/32/33 (content shifted
# down automatically by Insert). Their sequence number in column A needs to
# be bumped by two to stay in sync with the new running index.
$ws2.Cells.Item(31, 1).Value = 30
$ws2.Cells.Item(32, 1).Value = 31
$ws2.Cells.Item(33, 1).Value = 32

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F5").Value = 1661
$ws3.Range("F7").Value = 1002

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (All types) -- mirrors the same F-column bumps as the
# other three sheets (no row insert happens here).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 1661
$ws4.Range("F6").Value = 1002
$ws4.Range("F7").Value = 28
$ws4.Range("F8").Value = 799
$ws4.Range("F9").Value = 38
$ws4.Range("F12").Value = 1151
$ws4.Range("F13").Value = 908
$ws4.Range("F17").Value = 718
$ws4.Range("F18").Value = 174
$ws4.Range("F19").Value = 174
$ws4.Range("F22").Value = 1437
$ws4.Range("F23").Value = 56
$ws4.Range("F25").Value = 1603
$ws4.Range("F27").Value = 608
$ws4.Range("F30").Value = 1079
$ws4.Range("F31").Value = 1511
$ws4.Range("F32").Value = 752
$ws4.Range("F33").Value = 616
$ws4.Range("F34").Value = 489
$ws4.Range("F38").Value = 256
$ws4.Range("F39").Value = 257
$ws4.Range("F42").Value = 299
$ws4.Range("F43").Value = 2422
$ws4.Range("F48").Value = 1354
$ws4.Range("F49").Value = 460
$ws4.Range("F50").Value = 3949
